$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.583.34"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.573.85"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.991"
$ws.Range("E4").Value = "  -1.78%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.40"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.492"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.27"
$ws.Range("E8").Value = "  +6.03%  "
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0598"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.798.26"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.544.65"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.522"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.30"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.524.96"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.39"
$ws.Range("E18").Value = "  +7.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.53"
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0707"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("E21").Value = "  -1.56%  "
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.46"
$ws.Range("E23").Value = "  +3.13%  "
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.11"
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.37"
$ws.Range("E26").Value = "  +2.16%  "
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.992"
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0474"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.462.46"
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.13"
$ws.Range("E34").Value = "  -1.19%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.61"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").Value = "  -5.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.816"
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.38"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.68"
$ws.Range("E42").Value = "  -2.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.992"
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.84"
$ws.Range("E44").Value = "  +6.59%  "
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.29"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.708.50"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.92"
$ws.Range("E48").Value = "  +1.63%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0524"
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0999"
$ws.Range("E50").Value = "  -2.65%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "40.20"
$ws.Range("E51").Value = "  +18.88%  "
